$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" (D) and "is_enabled" (E) columns entirely.
# This shifts the old F (order_by) -> D and old G (rem) -> E.
$ws.Range("D1:E1").EntireColumn.Delete() | Out-Null

# The old "order_by" column (now in D) and "rem" column (now in E) stay as-is.
# Replace the last used column (now F, previously unused) with the new
# "tenant_id" column content, taking the place vacated by the deleted columns.
$ws.Range("F1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'
